# edit.ps1 - apply the "rapport fase 1" diff via Word COM-interop
# Runtime note: $word / $d are pre-seeded. wdFindContinue = 1, wdReplaceAll = 2.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
    return $found
}

$vt = [char]11   # vertical-tab == <w:br/> in Word's text stream

# 1) "Het systeem zal via het internet..." - merge the two runs (identical text,
#    just collapses the run split) by doing a no-op find/replace across both runs.
$s1 = "Het systeem zal via het internet bereikbaar zijn. Hierdoor is het makkelijk voor iedereen te gebruiken; er hoeft niet een speciaal programma geïnstalleerd te worden. Een nadeel hiervan zou wel kunnen zijn dat het niet mogelijk is om alle muziekbestanden die op de computer van de gebruiker zijn, door te scannen. "
Replace-Text $s1 $s1 | Out-Null

# 2) "een databron" -> "één data bron"
Replace-Text "Hiervoor zal gebruikt worden gemaakt van een databron. In latere fases" `
             "Hiervoor zal gebruikt worden gemaakt van één data bron. In latere fases" | Out-Null

# 3) Remove the 'Pop' sentence and join "...aangeven. " directly with "Uit deze selectie"
#    (also removes the <w:br/> that used to separate them).
$old3 = "), is een selectie gemaakt van tags welke vaak worden gebruikt en daarnaast ook nog een duidelijk genre aangeven. " + [char]8216 + "Pop" + [char]8217 + ", ofwel " + [char]8216 + "populaire muziek" + [char]8217 + " is bijvoorbeeld één van de tags waarvoor is gekozen dat deze niet aan de keuzelijst wordt toegevoegd. Het is weliswaar een veelgebruikte tag, maar ook een genre dat niet heel strak omlijnt is, maar sterk afhankelijk is van gebruikersinput en daarnaast ook sterk veranderlijk is. " + $vt + "Uit deze selectie"
$new3 = "), is een selectie gemaakt van tags welke vaak worden gebruikt en daarnaast ook nog een duidelijk genre aangeven. Uit deze selectie"
Replace-Text $old3 $new3 | Out-Null

# 4) "voor elk van de genres" -> "voor de meest favoriete 3 genres"
Replace-Text "wordt er voor elk van de genres bijbehorende artiesten opgezocht uit de MusicBrainz database. " `
             "wordt er voor de meest favoriete 3 genres bijbehorende artiesten opgezocht uit de MusicBrainz database. " | Out-Null

# 5) Append new sentences about MusicBrainz discarding queries.
Replace-Text ": de muziekdatabase verwerkt slechts één zoekopdracht per seconde." `
             ": de muziekdatabase verwerkt slechts één zoekopdracht per seconde. Daarnaast slaat MusicBrainz zoekopdrachten over, indien binnen een seconde weer een nieuwe aanroep wordt gedaan. Hierdoor gaat er veel data verloren. " | Out-Null

# 6) Merge runs around "Dit hoeft echter niet zo te zijn. Zo krijgt..." (no text change)
$s6 = " gekozen genre. Dit hoeft echter niet zo te zijn. Zo krijgt bijvoorbeeld een artiest uit het 1"
Replace-Text $s6 $s6 | Out-Null

# 7) Merge runs around "...waardering van 12. Hierdoor zou..." (no text change)
$s7 = " gekozen genre met waardering 4 een totale gewogen waardering van 12. Hierdoor zou de (beter gewaardeerde) artiest uit het 3"
Replace-Text $s7 $s7 | Out-Null

# 8) Split "waardering" into "waar" + "dering" around the relocated _GoBack bookmark.
$r8 = $d.Content
$found8 = $r8.Find.Execute("aflopende grootte van waar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found8) {
    $pos = $r8.End
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
} else {
    Write-Output "NOT FOUND: aflopende grootte van waar"
}

Write-Output "done"
